$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 5: Problem 4 / Part A -> Problem 4 / Part "-"
$t.Cell(5,2).Range.Text = "-"

# Row 6: Problem 4 / Part B -> Problem 5 / Part "-"
$t.Cell(6,1).Range.Text = "5"
$t.Cell(6,2).Range.Text = "-"

# Row 7: Problem 5 -> Problem 6
$t.Cell(7,1).Range.Text = "6"

# Row 8: Problem 6 -> Problem 7
$t.Cell(8,1).Range.Text = "7"

# Row 9: Problem 7 / Part A -> Problem 8 / Part "-"
$t.Cell(9,1).Range.Text = "8"
$t.Cell(9,2).Range.Text = "-"

# Row 10: Problem 7 / Part B -> Problem 9 / Part "-", plus solution text rewording
$t.Cell(10,1).Range.Text = "9"
$t.Cell(10,2).Range.Text = "-"

$d.Content.Find.Execute(
    "normality assumption is met. This means you need to do a qq-plot for each of the",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "normality assumption is met. This means you need to create a histogram for each",
    2) | Out-Null

$d.Content.Find.Execute(
    "groups. The qq-plots show that the groups are not perfectly normal, but they are",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "of the groups. The histograms show that the groups are not perfectly normal, but",
    2) | Out-Null

$d.Content.Find.Execute(
    "probably close enough to proceed with ANOVA.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "they are probably close enough to proceed with ANOVA.",
    2) | Out-Null

# Row 11: Problem 7 / Part C -> Problem 10 / Part "-"
$t.Cell(11,1).Range.Text = "10"
$t.Cell(11,2).Range.Text = "-"

# Row 12: Problem 7 / Part D -> Problem 11 / Part "-"
$t.Cell(12,1).Range.Text = "11"
$t.Cell(12,2).Range.Text = "-"

# Row 13: Problem 7 / Part E -> Problem 12 / Part "-"
$t.Cell(13,1).Range.Text = "12"
$t.Cell(13,2).Range.Text = "-"

# Row 14: Problem 7 / Part F -> Problem 13 / Part "-"
$t.Cell(14,1).Range.Text = "13"
$t.Cell(14,2).Range.Text = "-"

# Row 15: Problem 7 / Part G -> Problem 14 / Part "-"
$t.Cell(15,1).Range.Text = "14"
$t.Cell(15,2).Range.Text = "-"

# Row 16: Problem 8 / Part (empty) -> Problem 15 / Part "-"
$t.Cell(16,1).Range.Text = "15"
$t.Cell(16,2).Range.Text = "-"
$t.Cell(16,2).Range.ParagraphFormat.Style = "Compact"
$t.Cell(16,2).Range.ParagraphFormat.Alignment = "left"
